$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 1625
$ws.Range("I34").Value = 1625
$ws.Range("K34").Value = 1625
$ws.Range("M34").Value = -1422
$ws.Range("H36").Value = 1625
$ws.Range("I36").Value = 1625
$ws.Range("K36").Value = 1625
$ws.Range("M36").Value = -910
$ws.Range("H53").Value = 101.833336
$ws.Range("I53").Value = 88
$ws.Range("K53").Value = 88
$ws.Range("M53").Value = 549
$ws.Range("H58").Value = 550.8333
$ws.Range("I58").Value = 550.8333
$ws.Range("K58").Value = 1652.4999
$ws.Range("M58").Value = -1502.4999
$ws.Range("H76").Value = 900
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 900
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H92").Value = 2664.8333
$ws.Range("I92").Value = 1500
$ws.Range("J92").Value = 2897.8
$ws.Range("K92").Value = 1500
$ws.Range("L92").Value = 2897.8
$ws.Range("M92").Value = -252
$ws.Range("N92").Value = -5393.8
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H101").Value = 621.3333
$ws.Range("I101").Value = 444.5
$ws.Range("J101").Value = 975
$ws.Range("K101").Value = 1333.5
$ws.Range("L101").Value = 2925
$ws.Range("M101").Value = 288.5
$ws.Range("N101").Value = -6169
$ws.Range("H107").Value = 3728.8
$ws.Range("I107").Value = 3793.5
$ws.Range("K107").Value = 3793.5
$ws.Range("M107").Value = -1873.5
$ws.Range("H109").Value = 55000
$ws.Range("J109").Value = 55000
$ws.Range("L109").Value = 55000
$ws.Range("N109").Value = -57774
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H115").Value = 300
$ws.Range("I115").Value = 300
$ws.Range("K115").Value = 900
$ws.Range("M115").Value = 667
$ws.Range("H136").Value = 99995
$ws.Range("J136").Value = 99995
$ws.Range("L136").Value = 99995
$ws.Range("N136").Value = -110195
$ws.Range("H139").Value = 99995
$ws.Range("J139").Value = 99995
$ws.Range("L139").Value = 99995
$ws.Range("N139").Value = -110275

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2320.5
$ws.Range("I45").Value = 2320.5
$ws.Range("K45").Value = 2320.5
$ws.Range("M45").Value = -1943.5
$ws.Range("H62").Value = 10000
$ws.Range("J62").Value = 10000
$ws.Range("L62").Value = 10000
$ws.Range("N62").Value = -11248
$ws.Range("H65").Value = 10000
$ws.Range("J65").Value = 10000
$ws.Range("L65").Value = 30000
$ws.Range("N65").Value = -36240
$ws.Range("H74").Value = 3456
$ws.Range("I74").Value = 2124.6667
$ws.Range("J74").Value = 7450
$ws.Range("K74").Value = 2124.6667
$ws.Range("L74").Value = 7450
$ws.Range("M74").Value = -1250.6667
$ws.Range("N74").Value = -9198
$ws.Range("H77").Value = 3456
$ws.Range("I77").Value = 2124.6667
$ws.Range("J77").Value = 7450
$ws.Range("K77").Value = 10623.3335
$ws.Range("L77").Value = 37250
$ws.Range("M77").Value = -6255.333500000001
$ws.Range("N77").Value = -45986
$ws.Range("H110").Value = 994.5
$ws.Range("I110").Value = 999
$ws.Range("K110").Value = 999
$ws.Range("M110").Value = 1046
$ws.Range("H122").Value = 3585.3333
$ws.Range("I122").Value = 3628
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 10884
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -8434
$ws.Range("N122").Value = -15400

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 18835.5
$ws.Range("J88").Value = 18835.5
$ws.Range("L88").Value = 18835.5
$ws.Range("N88").Value = -19647.5
$ws.Range("H91").Value = 18835.5
$ws.Range("J91").Value = 18835.5
$ws.Range("L91").Value = 18835.5
$ws.Range("N91").Value = -21643.5
$ws.Range("H95").Value = 3750
$ws.Range("J95").Value = 3750
$ws.Range("L95").Value = 3750
$ws.Range("N95").Value = -9242
$ws.Range("H130").Value = 94997
$ws.Range("J130").Value = 94997
$ws.Range("L130").Value = 94997
$ws.Range("N130").Value = -105037

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 711.5714
$ws.Range("I16").Value = 663.5
$ws.Range("K16").Value = 663.5
$ws.Range("M16").Value = -376.5
$ws.Range("H113").Value = 711.5714
$ws.Range("I113").Value = 663.5
$ws.Range("K113").Value = 663.5
$ws.Range("M113").Value = 1506.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 7898.8
$ws.Range("J92").Value = 7898.8
$ws.Range("L92").Value = 7898.8
$ws.Range("N92").Value = -11642.8
$ws.Range("H113").Value = 1670
$ws.Range("I113").Value = 1670
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1670
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 500
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H56").Value = 16320.2
$ws.Range("J56").Value = 15250
$ws.Range("L56").Value = 15250
$ws.Range("N56").Value = -16632
$ws.Range("H105").Value = 5615
$ws.Range("J105").Value = 5615
$ws.Range("L105").Value = 5615
$ws.Range("N105").Value = -12603
$ws.Range("H127").Value = 54998
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 54998
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 54998
$ws.Range("M127").ClearContents()
$ws.Range("N127").Value = -64918
$ws.Range("H134").Value = 99995
$ws.Range("J134").Value = 99995
$ws.Range("L134").Value = 99995
$ws.Range("N134").Value = -110135
$ws.Range("H135").Value = 219998
$ws.Range("J135").Value = 219998
$ws.Range("L135").Value = 219998
$ws.Range("N135").Value = -230138

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 3551.6667
$ws.Range("J47").Value = 20000
$ws.Range("L47").Value = 20000
$ws.Range("N47").Value = -21144
$ws.Range("H51").Value = 29035
$ws.Range("I51").Value = 29035
$ws.Range("K51").Value = 29035
$ws.Range("M51").Value = -28525
$ws.Range("H58").Value = 30085
$ws.Range("I58").Value = 30085
$ws.Range("K58").Value = 30085
$ws.Range("M58").Value = -29777
$ws.Range("H94").Value = 26500
$ws.Range("J94").Value = 26500
$ws.Range("L94").Value = 26500
$ws.Range("N94").Value = -28302
$ws.Range("H100").Value = 38649.5
$ws.Range("I100").Value = 38649.5
$ws.Range("K100").Value = 77299
$ws.Range("M100").Value = -76758
$ws.Range("H133").Value = 13411666
$ws.Range("J133").Value = 13411666
$ws.Range("L133").Value = 13411666
$ws.Range("N133").Value = -13421786
